$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.664.22'
$ws.Range("E2").Value = '  -3.00%  '
$ws.Range("D3").Value = '1.741.53'
$ws.Range("E3").Value = '  -5.24%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '238.22'
$ws.Range("E5").Value = '  -8.48%  '
$ws.Range("E6").Value = '  +0.12%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5023'
$ws.Range("E7").Value = '  -5.72%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '41.67'
$ws.Range("E8").Value = '  -6.94%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2625'
$ws.Range("E9").Value = '  -12.81%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06102'
$ws.Range("E10").Value = '  -11.12%  '
$ws.Range("D11").Value = '1.748.29'
$ws.Range("E11").Value = '  -4.98%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.06956'
$ws.Range("E12").Value = '  -6.45%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '15.09'
$ws.Range("E13").Value = '  -13.42%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.473'
$ws.Range("E14").Value = '  -9.93%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5895'
$ws.Range("E15").Value = '  -19.67%  '
$ws.Range("E16").Value = '  -14.40%  '
$ws.Range("E17").Value = '  -0.04%  '
$ws.Range("E18").Value = '  +0.12%  '
$ws.Range("D19").Value = '25.677.48'
$ws.Range("E19").Value = '  -3.05%  '
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000006756'
$ws.Range("E20").Value = '  -14.55%  '
$ws.Range("B21").Value = 'Avalanche'
$ws.Range("C21").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.54'
$ws.Range("E21").Value = '  -17.01%  '
$ws.Range("D22").Value = '1.970.19'
$ws.Range("E22").Value = '  -5.49%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.024'
$ws.Range("E23").Value = '  -12.20%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.103'
$ws.Range("E24").Value = '  -12.35%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.068'
$ws.Range("E25").Value = '  -15.11%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '137.48'
$ws.Range("E26").Value = '  -3.52%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.533'
$ws.Range("E27").Value = '  -9.02%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.814'
$ws.Range("E28").Value = '  -18.21%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '14.94'
$ws.Range("E29").Value = '  -11.67%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '103.37'
$ws.Range("E30").Value = '  -6.25%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.753'
$ws.Range("E31").Value = '  -11.68%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08077'
$ws.Range("E32").Value = '  -7.97%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.443'
$ws.Range("E33").Value = '  -14.30%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04490'
$ws.Range("E34").Value = '  -6.12%  '
$ws.Range("E35").Value = '  +0.09%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.659'
$ws.Range("E36").Value = '  -9.03%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9734'
$ws.Range("E37").Value = '  -13.66%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6050'
$ws.Range("E38").Value = '  -17.06%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.646'
$ws.Range("E39").Value = '  -14.61%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01540'
$ws.Range("E40").Value = '  -9.93%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.914'
$ws.Range("E41").Value = '  -16.33%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.003'
$ws.Range("E42").Value = '  +0.08%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '103.67'
$ws.Range("E43").Value = '  -3.67%  '
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3777'
$ws.Range("E44").Value = '  -19.69%  '
$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.094'
$ws.Range("E45").Value = '  -13.17%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.7257'
$ws.Range("E46").Value = '  -20.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05332'
$ws.Range("E47").Value = '  -8.06%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1103'
$ws.Range("E48").Value = '  -10.35%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '29.91'
$ws.Range("E49").Value = '  -14.06%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.847'
$ws.Range("E50").Value = '  -20.37%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '52.22'
$ws.Range("E51").Value = '  -13.11%  '
